$d = $word.ActiveDocument

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

# ---------------------------------------------------------------------
# 1) Add a new "data licence acknowledgement" paragraph right after the
#    Gates Foundation acknowledgement paragraph (before the References
#    heading / the acknowledgements bookmark end).
# ---------------------------------------------------------------------
$ackPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Bill*Melinda Gates Foundation*") {
        $ackPara = $p
    }
}

if ($ackPara -eq $null) {
    throw "Could not locate the Gates Foundation acknowledgements paragraph"
}

$ackRange = $ackPara.Range
$ackRange.InsertParagraphAfter()

$newPara = $ackPara.Next()
$newPara.Style = "BodyText"

$newText = "The package associated with this paper contains information from the dataset " + $openQuote + "LTA MRT Station Exit (GEOJSON)" + $closeQuote + " accessed on the 10th of December 2024 from data.gov.sg, which is made available under the terms of the Singapore Open Data Licence version 1.0 https://data.gov.sg/open-data-licence."

$newPara.Range.Text = $newText

# ---------------------------------------------------------------------
# 2) Update the "LTA MRT Station Exit (GEOJSON) Dataset" bibliography
#    entry: italicise the title and append ". data.gov.sg." afterwards,
#    removing the surrounding curly quotes and trailing period.
# ---------------------------------------------------------------------
$searchText = $openQuote + "LTA MRT Station Exit (GEOJSON) Dataset." + $closeQuote

$refRange = $d.Content
$found = $refRange.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the LTA MRT Station Exit bibliography text"
}

$refStart = $refRange.Start

$italicText = "LTA MRT Station Exit (GEOJSON) Dataset"
$tailText = ". data.gov.sg."

$refRange.Text = $italicText + $tailText

$italicRange = $d.Range($refStart, $refStart + $italicText.Length)
$italicRange.Italic = $true

Write-Output "done"
